$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1300.4615
$ws.Cells.Item(17, 10).Value = 1300.4615
$ws.Cells.Item(17, 12).Value = 3901.3845
$ws.Cells.Item(17, 14).Value = -4237.3845
$ws.Cells.Item(106, 8).Value = 1848
$ws.Cells.Item(106, 9).Value = 1848
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 11).Value = 1848
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 13).ClearContents()
$ws.Cells.Item(106, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 37501960
$ws.Cells.Item(107, 9).Value = 17859584
$ws.Cells.Item(107, 11).Value = 17859584
$ws.Cells.Item(107, 13).Value = -17857664
$ws.Cells.Item(125, 8).Value = 125001544
$ws.Cells.Item(125, 9).Value = 200001060
$ws.Cells.Item(125, 10).Value = 2366.3333
$ws.Cells.Item(125, 11).Value = 1800009540
$ws.Cells.Item(125, 12).Value = 21296.9997
$ws.Cells.Item(125, 13).Value = -1800007080
$ws.Cells.Item(125, 14).Value = -26216.9997
$ws.Cells.Item(132, 8).Value = 1992.3889
$ws.Cells.Item(132, 9).Value = 1337.7693
$ws.Cells.Item(132, 11).Value = 4013.3079
$ws.Cells.Item(132, 13).Value = -1483.3079
$ws.Cells.Item(138, 8).Value = 1591182.6
$ws.Cells.Item(138, 10).Value = 3130618
$ws.Cells.Item(138, 12).Value = 9391854
$ws.Cells.Item(138, 14).Value = -9402134
$ws.Cells.Item(141, 8).Value = 4836.5
$ws.Cells.Item(141, 9).Value = 4352.4443
$ws.Cells.Item(141, 11).Value = 13057.3329
$ws.Cells.Item(141, 13).Value = -7877.332900000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4552719.5
$ws.Cells.Item(32, 9).Value = 5133685.5
$ws.Cells.Item(32, 10).Value = 21183.6
$ws.Cells.Item(32, 11).Value = 5133685.5
$ws.Cells.Item(32, 12).Value = 21183.6
$ws.Cells.Item(32, 13).Value = -5133398.5
$ws.Cells.Item(32, 14).Value = -21757.6
$ws.Cells.Item(122, 8).Value = 4621.905
$ws.Cells.Item(122, 9).Value = 2073.125
$ws.Cells.Item(122, 11).Value = 6219.375
$ws.Cells.Item(122, 13).Value = -3769.375
$ws.Cells.Item(139, 8).Value = 65916.25
$ws.Cells.Item(139, 10).Value = 65916.25
$ws.Cells.Item(139, 12).Value = 65916.25
$ws.Cells.Item(139, 14).Value = -76196.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 48617204
$ws.Cells.Item(20, 9).Value = 58336944
$ws.Cells.Item(20, 11).Value = 58336944
$ws.Cells.Item(20, 13).Value = -58336697
$ws.Cells.Item(134, 8).Value = 4799.2764
$ws.Cells.Item(134, 9).Value = 1745.9688
$ws.Cells.Item(134, 11).Value = 5237.9064
$ws.Cells.Item(134, 13).Value = -2702.9064

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 5832.1934
$ws.Cells.Item(16, 9).Value = 768.6
$ws.Cells.Item(16, 11).Value = 768.6
$ws.Cells.Item(16, 13).Value = -481.6
$ws.Cells.Item(50, 8).Value = 90089.664
$ws.Cells.Item(50, 10).Value = 90089.664
$ws.Cells.Item(50, 12).Value = 90089.664
$ws.Cells.Item(50, 14).Value = -91339.664
$ws.Cells.Item(107, 8).Value = 1113.9333
$ws.Cells.Item(107, 9).Value = 326.55554
$ws.Cells.Item(107, 10).Value = 2295
$ws.Cells.Item(107, 11).Value = 326.55554
$ws.Cells.Item(107, 12).Value = 2295
$ws.Cells.Item(107, 13).Value = 1593.44446
$ws.Cells.Item(107, 14).Value = -6135
$ws.Cells.Item(113, 8).Value = 5832.1934
$ws.Cells.Item(113, 9).Value = 768.6
$ws.Cells.Item(113, 11).Value = 768.6
$ws.Cells.Item(113, 13).Value = 1401.4
$ws.Cells.Item(122, 8).Value = 2205.625
$ws.Cells.Item(122, 9).Value = 1514.6
$ws.Cells.Item(122, 11).Value = 4543.799999999999
$ws.Cells.Item(122, 13).Value = -2093.799999999999
$ws.Cells.Item(132, 8).Value = 4739.0225
$ws.Cells.Item(132, 9).Value = 2169.9583
$ws.Cells.Item(132, 11).Value = 6509.874899999999
$ws.Cells.Item(132, 13).Value = -3979.874899999999
$ws.Cells.Item(134, 8).Value = 5028.514
$ws.Cells.Item(134, 9).Value = 1431.0454
$ws.Cells.Item(134, 11).Value = 4293.1362
$ws.Cells.Item(134, 13).Value = -1758.1362

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(81, 8).Value = 188248.5
$ws.Cells.Item(81, 10).Value = 188248.5
$ws.Cells.Item(81, 12).Value = 564745.5
$ws.Cells.Item(81, 14).Value = -566991.5
$ws.Cells.Item(84, 8).Value = 188248.5
$ws.Cells.Item(84, 10).Value = 188248.5
$ws.Cells.Item(84, 12).Value = 1694236.5
$ws.Cells.Item(84, 14).Value = -1705468.5
$ws.Cells.Item(107, 8).Value = 881.2941
$ws.Cells.Item(107, 9).Value = 675
$ws.Cells.Item(107, 10).Value = 925.5
$ws.Cells.Item(107, 11).Value = 2025
$ws.Cells.Item(107, 12).Value = 2776.5
$ws.Cells.Item(107, 13).Value = -105
$ws.Cells.Item(107, 14).Value = -6616.5
$ws.Cells.Item(119, 8).Value = 1809.3334
$ws.Cells.Item(119, 9).Value = 214
$ws.Cells.Item(119, 11).Value = 642
$ws.Cells.Item(119, 13).Value = 4196
$ws.Cells.Item(120, 8).Value = 28479.75
$ws.Cells.Item(120, 9).Value = 5306.6665
$ws.Cells.Item(120, 10).Value = 97999
$ws.Cells.Item(120, 11).Value = 15919.9995
$ws.Cells.Item(120, 12).Value = 293997
$ws.Cells.Item(120, 13).Value = -11081.9995
$ws.Cells.Item(120, 14).Value = -303673
$ws.Cells.Item(132, 8).Value = 4722.2793
$ws.Cells.Item(132, 10).Value = 6235.407
$ws.Cells.Item(132, 12).Value = 56118.663
$ws.Cells.Item(132, 14).Value = -61178.663

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 4291.5557
$ws.Cells.Item(102, 9).Value = 4072.261
$ws.Cells.Item(102, 10).Value = 5552.5
$ws.Cells.Item(102, 11).Value = 4072.261
$ws.Cells.Item(102, 12).Value = 5552.5
$ws.Cells.Item(102, 13).Value = -2450.261
$ws.Cells.Item(102, 14).Value = -8796.5
$ws.Cells.Item(113, 8).Value = 7303.788
$ws.Cells.Item(113, 10).Value = 9966.294
$ws.Cells.Item(113, 12).Value = 9966.294
$ws.Cells.Item(113, 14).Value = -14306.294
$ws.Cells.Item(122, 8).Value = 82638.92999999999
$ws.Cells.Item(122, 9).Value = 153062
$ws.Cells.Item(122, 10).Value = 2155.4285
$ws.Cells.Item(122, 11).Value = 459186
$ws.Cells.Item(122, 12).Value = 6466.2855
$ws.Cells.Item(122, 13).Value = -456736
$ws.Cells.Item(122, 14).Value = -11366.2855
$ws.Cells.Item(126, 8).Value = 4410.095
$ws.Cells.Item(126, 10).Value = 6481.4
$ws.Cells.Item(126, 12).Value = 19444.2
$ws.Cells.Item(126, 14).Value = -24384.2
$ws.Cells.Item(132, 8).Value = 4195.08
$ws.Cells.Item(132, 9).Value = 1898.9048
$ws.Cells.Item(132, 10).Value = 16250
$ws.Cells.Item(132, 11).Value = 5696.7144
$ws.Cells.Item(132, 12).Value = 48750
$ws.Cells.Item(132, 13).Value = -3166.7144
$ws.Cells.Item(132, 14).Value = -53810

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5106.0835
$ws.Cells.Item(7, 9).Value = 2939.8572
$ws.Cells.Item(7, 11).Value = 2939.8572
$ws.Cells.Item(7, 13).Value = -2827.8572
$ws.Cells.Item(16, 8).Value = 954.8889
$ws.Cells.Item(16, 9).Value = 956.2857
$ws.Cells.Item(16, 11).Value = 956.2857
$ws.Cells.Item(16, 13).Value = -786.2857
$ws.Cells.Item(46, 8).Value = 2999.3044
$ws.Cells.Item(46, 9).Value = 2488.9333
$ws.Cells.Item(46, 10).Value = 3956.25
$ws.Cells.Item(46, 11).Value = 2488.9333
$ws.Cells.Item(46, 12).Value = 3956.25
$ws.Cells.Item(46, 13).Value = -2300.9333
$ws.Cells.Item(46, 14).Value = -4332.25
$ws.Cells.Item(61, 8).Value = 5025.55
$ws.Cells.Item(61, 9).Value = 3729.25
$ws.Cells.Item(61, 11).Value = 3729.25
$ws.Cells.Item(61, 13).Value = -3527.25
$ws.Cells.Item(113, 8).Value = 5025.55
$ws.Cells.Item(113, 9).Value = 3729.25
$ws.Cells.Item(113, 11).Value = 3729.25
$ws.Cells.Item(113, 13).Value = -1559.25
$ws.Cells.Item(122, 8).Value = 4505.3
$ws.Cells.Item(122, 9).Value = 2720.25
$ws.Cells.Item(122, 10).Value = 5695.3335
$ws.Cells.Item(122, 11).Value = 8160.75
$ws.Cells.Item(122, 12).Value = 17086.0005
$ws.Cells.Item(122, 13).Value = -5710.75
$ws.Cells.Item(122, 14).Value = -21986.0005
$ws.Cells.Item(126, 8).Value = 5106.0835
$ws.Cells.Item(126, 9).Value = 2939.8572
$ws.Cells.Item(126, 11).Value = 8819.571599999999
$ws.Cells.Item(126, 13).Value = -6349.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(74, 8).Value = 5711.6
$ws.Cells.Item(74, 10).Value = 5711.6
$ws.Cells.Item(74, 12).Value = 5711.6
$ws.Cells.Item(74, 14).Value = -7583.6
$ws.Cells.Item(77, 8).Value = 5711.6
$ws.Cells.Item(77, 10).Value = 5711.6
$ws.Cells.Item(77, 12).Value = 17134.8
$ws.Cells.Item(77, 14).Value = -26494.8
$ws.Cells.Item(107, 8).Value = 11495300
$ws.Cells.Item(107, 10).Value = 22223750
$ws.Cells.Item(107, 12).Value = 66671250
$ws.Cells.Item(107, 14).Value = -66675090
$ws.Cells.Item(122, 8).Value = 3033.9412
$ws.Cells.Item(122, 9).Value = 2352.875
$ws.Cells.Item(122, 10).Value = 4668.5
$ws.Cells.Item(122, 11).Value = 7058.625
$ws.Cells.Item(122, 12).Value = 14005.5
$ws.Cells.Item(122, 13).Value = -4608.625
$ws.Cells.Item(122, 14).Value = -18905.5
$ws.Cells.Item(136, 8).Value = 3816.2888
$ws.Cells.Item(136, 9).Value = 1432.5834
$ws.Cells.Item(136, 11).Value = 4297.7502
$ws.Cells.Item(136, 13).Value = -1747.7502
